$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title
Replace-Text "Unraveling the Enigmatic Quantum Realm" "Exploring the Realm of Chemistry: A Journey Through the Elements and Beyond"

# Author name (3 runs -> merges into a single run "Emily Rodriguez")
Replace-Text "Dr. Eleanor Abernathy" "Emily Rodriguez"

# Email user + domain
Replace-Text "eabernathy@quantumresearch" "rodriguez_e@xyz"
Replace-Text "org" "edu"

# Body paragraph 1
Replace-Text "In the captivating realm of quantum mechanics, the universe unveils a symphony of enigmatic phenomena that challenge our conventional understanding of reality" "In the vast realm of science, Chemistry unravels the intricate world of elements, their interactions, and the transformation of substances"

Replace-Text "As we delve into the subatomic world, we encounter particles that exist in multiple states simultaneously, perplexing entanglement, and the mysterious influence of the observer" "It probes the fundamental principles behind the composition, structure, and properties of matter"

Replace-Text "These mind-bending concepts, defying classical logic, have profoundly impacted our comprehension of the universe, ushering in a new era of scientific exploration and innovation" "Unraveling the secrets of chemistry allows us to unravel the secrets of our vast world. As we embark on this captivating journey through the elements and beyond, we will discover how chemistry shapes our everyday experiences and unlocks the potential for unimaginable discoveries"

# Body paragraph 2
Replace-Text "The advent of quantum theory has illuminated the intricate ballet of particles, revealing a universe teeming with uncertainty and probability" "Chemistry permeates our lives in countless ways"

Replace-Text "Quantum mechanics has shattered the illusion of solid, well-defined particles, replacing it with a hazy realm of probabilities and wave functions" "From the food we consume to the clothes we wear, the medicines that heal us, and the fuels that power our vehicles, chemistry is the foundation of modern society"

Replace-Text "This fundamental shift in our understanding of matter has opened up unprecedented possibilities for technological breakthroughs, from quantum computing to ultra-precise sensors" "Understanding the principles of chemistry enables us to appreciate the profound impact it has on our daily lives and the ways it influences the world around us. Whether we are unraveling the complexities of a chemical reaction or exploring the vastness of the periodic table, chemistry invites us to explore the mysteries of the natural world"

# Body paragraph 3
Replace-Text "Moreover, quantum mechanics has sparked a profound reexamination of consciousness and the nature of reality" "The study of chemistry fosters a sense of curiosity, critical thinking, and problem-solving skills, nurturing our understanding of the universe and inspiring countless innovations"

Replace-Text "The enigmatic nature of quantum phenomena, such as superposition and entanglement, has ignited debates among physicists, philosophers, and theologians, leading to new perspectives on the fundamental questions of existence and consciousness" "From the development of new materials and technologies to the pursuit of sustainable solutions, chemistry holds the key to addressing some of the world's most pressing challenges"

Replace-Text "The study of quantum mechanics has not only revolutionized our understanding of the physical world but has also challenged our most deeply held assumptions about reality" "By embracing the incredible world of chemistry, we unlock the potential to shape a brighter future for ourselves and generations to come"

# Summary heading paragraph content
Replace-Text "Quantum mechanics, with its enigmatic phenomena and mind-bending concepts, has transformed our understanding of the universe" "Chemistry, revealing the intricacies of matter and its transformations, is integral to unraveling the secrets of our world and addressing pressing challenges"

Replace-Text "From the perplexing world of particles to the profound implications for consciousness and reality, quantum theory continues to inspire scientific exploration and philosophical contemplation" "Exploring the principles of chemistry exposes us to the foundation of modern society, changing how we comprehend and interact with the world around us"

Replace-Text "Its impact knows no bounds, extending from technological advancements to fundamental questions of existence. As we venture further into the enigmatic quantum realm, we unravel the secrets of the subatomic world, forever changing our perception of reality" "Through the study of chemistry, we can appreciate nature's exquisite tapestry, nurture our curiosity, hone our problem-solving skills, and work towards shaping a brighter, more sustainable future"

# Add a new empty paragraph at the very end of the document body
$end = $d.Content
$end.Collapse(0)
$end.InsertParagraphAfter()

Write-Output "edit complete"
